$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Periodo Mora" column (E16:E28): the table of worker dues periods is
# reordered from descending (2103 .. 2003) to ascending (2003 .. 2103).
$ws.Range("E16").Value = "2003"
$ws.Range("E17").Value = "2004"
$ws.Range("E18").Value = "2005"
$ws.Range("E19").Value = "2006"
$ws.Range("E20").Value = "2007"
$ws.Range("E21").Value = "2008"
$ws.Range("E22").Value = "2009"
$ws.Range("E23").Value = "2010"
$ws.Range("E24").Value = "2011"
$ws.Range("E25").Value = "2012"
$ws.Range("E26").Value = "2101"
$ws.Range("E27").Value = "2102"
$ws.Range("E28").Value = "2103"

# "Valor Mora" (column F) values follow the same reordering: the amount that
# used to belong to period 2103 (row 16) now belongs to period 2003 (still
# row 16), and the amount that used to belong to period 2003 (row 28) now
# belongs to period 2103 (still row 28).
$ws.Range("F16").Value = 35112
$ws.Range("F28").Value = 26919
